$d = $word.ActiveDocument

# --- 1. Mark the three screenshot picture runs as NoProof (adds <w:noProof/> to rPr) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $para.Range.NoProofing = 1
    }
}

# --- 2. Locate the "URL to GitHub Repository:" paragraph ---
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd() -eq "URL to GitHub Repository:") {
        $targetIndex = $i
    }
}
$target = $d.Paragraphs($targetIndex)

# --- 3. Remove the existing _GoBack bookmark; it currently sits inside the last
#        picture's paragraph and needs to move to the end of the new URL text ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 4. Insert a new paragraph after it containing the repository URL. An extra
#        trailing placeholder character is typed along with the real text in a single
#        Range.Text assignment (so run formatting - <w:b/><w:szCs w:val="24"/> - is
#        inherited correctly from the paragraph), a bookmark is anchored just before
#        that placeholder, and the placeholder is deleted afterwards - the bookmark's
#        position shifts (gravitates) left with it and survives exactly at the end of
#        the real URL text. ---
$afterInsertPos = $target.Range.End
$target.Range.InsertParagraphAfter()

$urlParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $afterInsertPos) {
        $urlParaIndex = $i
        break
    }
}

$urlPara = $d.Paragraphs($urlParaIndex)
$urlPara.Range.Text = "https://github.com/mctimoth/FESD-FET-Week1#"

$urlPara = $d.Paragraphs($urlParaIndex)
$bookmarkPos = $urlPara.Range.End - 2
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$urlPara = $d.Paragraphs($urlParaIndex)
$placeholderPos = $urlPara.Range.End - 2
$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()

Write-Output "done"
